# DBDT.xlsx - "CPME Acetone Water" sheet, row 11:
# Fill in the measured mass-fraction inputs (D11:I11) and the computed
# mole-fraction / normalized-fraction formulas (J11:U11), mirroring the
# pattern already present in rows above (e.g. row 10), then move the
# active selection to J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raw measured values (previously blank)
$ws.Range("D11").Value = 0.22484456
$ws.Range("E11").Value = 0.37722371999999998
$ws.Range("F11").Value = 0.39793171999999999
$ws.Range("G11").Value = 0.01691051
$ws.Range("H11").Value = 0.17791449000000001
$ws.Range("I11").Value = 0.80517499000000003

# Derived ratios (previously plain pasted values, now live formulas)
$ws.Range("J11").Formula = '=(P11/$A$6)/((P11/$A$6)+(Q11/$B$6)+(R11/$C$6))'
$ws.Range("K11").Formula = '=(Q11/$B$6)/((P11/$A$6)+(Q11/$B$6)+(R11/$C$6))'
$ws.Range("L11").Formula = '=(R11/$C$6)/((P11/$A$6)+(Q11/$B$6)+(R11/$C$6))'
$ws.Range("M11").Formula = '=(S11/$A$6)/((S11/$A$6)+(T11/$B$6)+(U11/$C$6))'
$ws.Range("N11").Formula = '=(T11/$B$6)/((S11/$A$6)+(T11/$B$6)+(U11/$C$6))'
$ws.Range("O11").Formula = '=(U11/$C$6)/((S11/$A$6)+(T11/$B$6)+(U11/$C$6))'

# Mole-fraction helper columns (previously blank)
$ws.Range("P11").Formula = '=(D11*$A$4)/((D11*$A$4)+(E11*$B$4)+(F11*$C$4))'
$ws.Range("Q11").Formula = '=(E11*$B$4)/((D11*$A$4)+(E11*$B$4)+(F11*$C$4))'
$ws.Range("R11").Formula = '=(F11*$C$4)/((D11*$A$4)+(E11*$B$4)+(F11*$C$4))'
$ws.Range("S11").Formula = '=(G11*$A$4)/((G11*$A$4)+(H11*$B$4)+(I11*$C$4))'
$ws.Range("T11").Formula = '=(H11*$B$4)/((G11*$A$4)+(H11*$B$4)+(I11*$C$4))'
$ws.Range("U11").Formula = '=(I11*$C$4)/((G11*$A$4)+(H11*$B$4)+(I11*$C$4))'

# Move/collapse the selection to J11 (was L14:L15 with active cell L15)
$ws.Range("J11").Select()
